# Adds the new "Core" attribute / miscellaneous translations to the DLP
# woordenlijst table (Tabel1), growing it from 14 data rows (A1:C15) to
# 23 data rows (A1:C24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# English term, Dutch term (uitleg remarks are added afterwards, below)
$newRows = @(
    @("Stewardship", "Rentmeesterschap"),
    @("Diplomacy",   "Diplomatie"),
    @("Martial",     "Krijgshaftigheid"),
    @("Intrigue",    "Intrige"),
    @("Learning",    "Geleerdheid"),
    @("Prowess",     "Dapperheid"),
    @("Dismiss",     "Ontzet"),
    @("Reject",      "Afwijzen"),
    @("Aptitude",    "Aanleg")
)

# First pass: add rows 1-8 (Stewardship .. Reject) with just their
# English/Dutch term, matching the order the terms were authored.
$addedRows = @()
for ($i = 0; $i -lt $newRows.Count - 1; $i++) {
    $rowData = $newRows[$i]
    $newRow = $tbl.ListRows.Add()
    $newRow.Range.Cells.Item(1, 1).Value = $rowData[0]
    $newRow.Range.Cells.Item(1, 2).Value = $rowData[1]
    $addedRows += , $newRow
}

# Second pass: fill in the "uitleg" cross-reference notes for the
# Dismiss/Reject pair, added after all the term pairs above.
$addedRows[6].Range.Cells.Item(1, 3).Value = "Voor Afwijzen zie: Reject"
$addedRows[7].Range.Cells.Item(1, 3).Value = "Voor Ontzet zie: Dismiss"

# Finally add the last row (Aptitude / Aanleg).
$lastData = $newRows[$newRows.Count - 1]
$finalRow = $tbl.ListRows.Add()
$finalRow.Range.Cells.Item(1, 1).Value = $lastData[0]
$finalRow.Range.Cells.Item(1, 2).Value = $lastData[1]

# Match the author's final selection position (one row below the last entry)
$selRow = $tbl.Range.Row + $tbl.Range.Rows.Count
$ws.Cells.Item($selRow, 1).Select()
